# "Got the barplots working with midas taxonomy"
#
# On the "Sheet4" worksheet, flag rows 102-118 in column J ("SCFA sample
# made") with the same "x" checkmark / green-fill style already used
# throughout the rest of that column (e.g. J2), and leave the selection
# on the cell the author ended up at (O112).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet4")

# Use an existing, correctly-formatted "x" flag cell in column J as the
# formatting template so the new cells pick up the same fill/alignment
# (style index 6 in the saved file) instead of Excel's bare default style.
$ws.Range("J2").Copy()

102..118 | ForEach-Object {
    $cell = $ws.Cells.Item($_, 10)   # column J = 10
    $cell.Value = "x"
    $cell.PasteSpecial(-4122)        # xlPasteFormats
}

$excel.CutCopyMode = $false

# Reflect the author's final selection after making the edits.
$ws.Range("O112").Select()
